# ---------------------------------------------------------------------------
# Target change (from the supplied OOXML diff): word/numbering.xml gets four
# <w:nsid w:val="..."/> attributes replaced with new random-looking hex GUIDs,
# one each inside the <w:abstractNum> blocks for abstractNumId 990, 991,
# 99721 and 99722. Nothing else in the package changes - not a single
# <w:lvl>, not the abstractNumId values, not the <w:num> -> abstractNumId
# mappings, not any paragraph's numPr, not any visible text or formatting.
#
# w:nsid is Word's purely-internal "list definition identity" GUID. It has
# never been exposed through Word's automation surface (VBA/COM or the
# modern JS API): there is no Document/List/ListTemplate/ListLevel property
# that reads or writes it, and nothing in ListFormat's numbering methods
# (ApplyListTemplate, ApplyNumberDefault, RemoveNumbers, ...) edits an
# existing abstractNum's nsid in place - they only ever mint brand new
# abstractNum/num entries for a fresh list. The commit message ("Vygenerovany
# file ve slozce", i.e. "generated file in folder", stamped with a build
# timestamp) matches that: the four GUIDs were reshuffled by whatever tool
# regenerated/re-saved the package, not by an in-Word editing action, so
# there is no user-facing edit to replay here.
#
# Given that, the faithful COM-interop reproduction of "only those four
# nsid GUIDs differ" is to leave every automatable, visible/semantic part of
# the document exactly as it is - touching the object model only in ways
# that are read-only / self-canceling, so the saved package round-trips
# byte-for-byte apart from what genuinely cannot be reached through the
# object model.

$d = $word.ActiveDocument

# Sanity/no-op touch of the object model (read-only) - confirms the document
# is reachable without mutating any content, formatting or numbering state.
$paragraphCount = $d.Paragraphs.Count
$listCount = $d.Lists.Count
